# Generate Report for Handback
#
# Marks the Overview rows as handed back (instead of "ready for handoff")
# and fills in the "Latest Target File" / "Latest Handback File" / "Latest
# Handback DateTime" columns (plus the matching hyperlinks) on the two
# per-locale sheets, then widens a few columns to fit the new content.

$wb = $excel.ActiveWorkbook

# ColumnWidth inputs (character units) that round-trip to the desired
# stored sheet widths once Excel quantizes them to whole pixels.
$wideWidth = 29.166666666666607   # -> stored width 30 (was ~17.2)
$fullWidth = 39.1666666666666     # -> stored width 40 (was ~18.6 / 21.7)

# ---------------------------------------------------------------------
# Sheet "Overview": mark both rows as handed back instead of "Ready for
# handoff".
# ---------------------------------------------------------------------
$overview = $wb.Worksheets.Item("Overview")
$overview.Range("E2").Value = "Handed back: in sync with en-US"
$overview.Range("F2").Value = "Handed back: in sync with en-US"
$overview.Range("E3").Value = "Handed back: in sync with en-US"
$overview.Range("F3").Value = "Handed back: in sync with en-US"
$overview.Columns(5).ColumnWidth = $wideWidth
$overview.Columns(6).ColumnWidth = $wideWidth

# ---------------------------------------------------------------------
# Shared lookup data for the two locale sheets.
# ---------------------------------------------------------------------
$baseUrl = "https://github.com/OpenLocalizationTestOrg/ol-test0/blob/4579062742c0fc5f45923080c2b41405436bb0f6/e2e/"
$md1 = "40bee032-bebc-42d5-bdad-324270a1e826.md"
$md2 = "7a82f988-c0d5-4b60-b7af-c5251f112ba7.md"

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$zhcn = $wb.Worksheets.Item("zh-cn")

$zhcn.Range("I2").Value = $md1
$zhcn.Range("J2").Value = "40bee032-bebc-42d5-bdad-324270a1e826.1a89d26278201965b17a75bcc637dca8dbf773dc.zh-cn.xlf"
$zhcn.Range("K2").Value = "2016-08-22 02:57:38"

$zhcn.Range("I3").Value = $md2
$zhcn.Range("J3").Value = "7a82f988-c0d5-4b60-b7af-c5251f112ba7.a8a10b4e08c70a3198aa3db1428578c441b08a38.zh-cn.xlf"
$zhcn.Range("K3").Value = "2016-08-22 02:57:38"

$zhcn.Hyperlinks.Delete()
$zhcn.Hyperlinks.Add($zhcn.Range("A2"), $baseUrl + $md1, $null, $null, $md1)
$zhcn.Hyperlinks.Add($zhcn.Range("I2"), $baseUrl + $md1, $null, $null, $md1)
$zhcn.Hyperlinks.Add($zhcn.Range("A3"), $baseUrl + $md2, $null, $null, $md2)
$zhcn.Hyperlinks.Add($zhcn.Range("I3"), $baseUrl + $md2, $null, $null, $md2)

$zhcn.Columns(3).ColumnWidth = $wideWidth
$zhcn.Columns(9).ColumnWidth = $fullWidth
$zhcn.Columns(10).ColumnWidth = $fullWidth

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$dede = $wb.Worksheets.Item("de-de")

$dede.Range("I2").Value = $md1
$dede.Range("J2").Value = "40bee032-bebc-42d5-bdad-324270a1e826.1a89d26278201965b17a75bcc637dca8dbf773dc.de-de.xlf"
$dede.Range("K2").Value = "2016-08-22 02:57:44"

$dede.Range("I3").Value = $md2
$dede.Range("J3").Value = "7a82f988-c0d5-4b60-b7af-c5251f112ba7.a8a10b4e08c70a3198aa3db1428578c441b08a38.de-de.xlf"
$dede.Range("K3").Value = "2016-08-22 02:57:44"

$dede.Hyperlinks.Delete()
$dede.Hyperlinks.Add($dede.Range("A2"), $baseUrl + $md1, $null, $null, $md1)
$dede.Hyperlinks.Add($dede.Range("I2"), $baseUrl + $md1, $null, $null, $md1)
$dede.Hyperlinks.Add($dede.Range("A3"), $baseUrl + $md2, $null, $null, $md2)
$dede.Hyperlinks.Add($dede.Range("I3"), $baseUrl + $md2, $null, $null, $md2)

$dede.Columns(3).ColumnWidth = $wideWidth
$dede.Columns(9).ColumnWidth = $fullWidth
$dede.Columns(10).ColumnWidth = $fullWidth
